# True Value CR for Resources Changes Added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add two new columns ---
$ws.Range("O1").Value = "ResourcesInAppPath"
$ws.Range("P1").Value = "ResourcesOutsidePath"

# --- Data row (row 2): update existing values, left to right ---
$ws.Range("A2").Value = "Chrome"
$ws.Range("B2").Value = "https://cmstest.webspiders.com/"
$ws.Range("F2").Value = "Ethan.taylor@yopmail.com"
$ws.Range("H2").Value = "SGP Summit"
$ws.Range("I2").Value = "SGPS18"
$ws.Range("J2").Value = "Registration"

# --- Data row (row 2): new columns ---
$ws.Range("O2").Value = "/Test Data/CustomResourceInApp.xlsx"
$ws.Range("P2").Value = "/Test Data/CustomResourceOutSide.xlsx"

# --- Formatting for the newly added header cells (bold header font, no border) ---
$ws.Range("O1").Font.Bold = $true
$ws.Range("O1").Font.Size = 12
$ws.Range("O1").Font.Name = "Arial"
$ws.Range("P1").Font.Bold = $true
$ws.Range("P1").Font.Size = 12
$ws.Range("P1").Font.Name = "Arial"

# --- Formatting for the newly added data cells: match the rest of row 2 (copy from N2) ---
$ws.Range("N2").Copy()
$ws.Range("O2:P2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Hyperlinks ---
# Update the existing B2 hyperlink's display text to match the new URL value.
$ws.Hyperlinks.Item(1).TextToDisplay = "https://cmstest.webspiders.com/"

# Add a new hyperlink on the email cell, then restore its original (non-hyperlink) look.
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:Ethan.taylor@yopmail.com", "", "", "Ethan.taylor@yopmail.com")
$ws.Range("F2").Font.Name = "Arial"
$ws.Range("F2").Font.Size = 11
$ws.Range("F2").Font.Underline = $false
$ws.Range("F2").Font.Color = 16711722

# --- Move selection to the new last cell, matching the recorded view state ---
$ws.Range("P1").Select()
